$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-01-08 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-14 Tuesday", 2) | Out-Null

# Update the division problems in the table, cell by cell (row, col are 1-based)
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Text = "63÷3=21, 0"

$cell = $t.Cell(1, 2)
$cell.Range.Text = "30÷6=5, 0"

$cell = $t.Cell(1, 3)
$cell.Range.Text = "89÷5=17, 4"

$cell = $t.Cell(1, 4)
$cell.Range.Text = "72÷4=18, 0"

$cell = $t.Cell(1, 5)
$cell.Range.Text = "72÷6=12, 0"

$cell = $t.Cell(5, 1)
$cell.Range.Text = "72÷7=10, 2"

$cell = $t.Cell(5, 2)
$cell.Range.Text = "22÷6=3, 4"

$cell = $t.Cell(5, 3)
$cell.Range.Text = "75÷5=15, 0"

$cell = $t.Cell(5, 4)
$cell.Range.Text = "39÷6=6, 3"

$cell = $t.Cell(5, 5)
$cell.Range.Text = "70÷3=23, 1"

$cell = $t.Cell(9, 1)
$cell.Range.Text = "24÷2=12, 0"

$cell = $t.Cell(9, 2)
$cell.Range.Text = "97÷4=24, 1"

$cell = $t.Cell(9, 3)
$cell.Range.Text = "55÷8=6, 7"

$cell = $t.Cell(9, 4)
$cell.Range.Text = "65÷7=9, 2"

$cell = $t.Cell(9, 5)
$cell.Range.Text = "71÷5=14, 1"

$cell = $t.Cell(13, 1)
$cell.Range.Text = "88÷9=9, 7"

$cell = $t.Cell(13, 2)
$cell.Range.Text = "81÷3=27, 0"

$cell = $t.Cell(13, 3)
$cell.Range.Text = "69÷8=8, 5"

$cell = $t.Cell(13, 4)
$cell.Range.Text = "12÷8=1, 4"

$cell = $t.Cell(13, 5)
$cell.Range.Text = "74÷8=9, 2"

$cell = $t.Cell(17, 1)
$cell.Range.Text = "75÷9=8, 3"

$cell = $t.Cell(17, 2)
$cell.Range.Text = "43÷7=6, 1"

$cell = $t.Cell(17, 3)
$cell.Range.Text = "21÷4=5, 1"

$cell = $t.Cell(17, 4)
$cell.Range.Text = "88÷7=12, 4"

$cell = $t.Cell(17, 5)
$cell.Range.Text = "90÷6=15, 0"

